$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 11024.211
$ws.Range("I70").Value = 2274.5557
$ws.Range("J70").Value = 18898.9
$ws.Range("K70").Value = 6823.6671
$ws.Range("L70").Value = 56696.7
$ws.Range("M70").Value = -6553.6671
$ws.Range("N70").Value = -57236.7
$ws.Range("H73").Value = 11024.211
$ws.Range("I73").Value = 2274.5557
$ws.Range("J73").Value = 18898.9
$ws.Range("K73").Value = 6823.6671
$ws.Range("L73").Value = 56696.7
$ws.Range("M73").Value = -5887.6671
$ws.Range("N73").Value = -58568.7
$ws.Range("H100").Value = 5267
$ws.Range("I100").Value = 2633.875
$ws.Range("K100").Value = 2633.875
$ws.Range("M100").Value = -2092.875
$ws.Range("H116").Value = 4363.778
$ws.Range("I116").Value = 4132.8125
$ws.Range("K116").Value = 4132.8125
$ws.Range("M116").Value = -690.8125
$ws.Range("H137").Value = 126224
$ws.Range("J137").Value = 5202
$ws.Range("L137").Value = 15606
$ws.Range("N137").Value = -20706
$ws.Range("H138").Value = 2819.6865
$ws.Range("I138").Value = 1298.15
$ws.Range("J138").Value = 3467.149
$ws.Range("K138").Value = 3894.45
$ws.Range("L138").Value = 10401.447
$ws.Range("M138").Value = 1245.55
$ws.Range("N138").Value = -20681.447

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10594.739
$ws.Range("I2").Value = 6860.9
$ws.Range("K2").Value = 6860.9
$ws.Range("M2").Value = -6747.9
$ws.Range("H116").Value = 10594.739
$ws.Range("I116").Value = 6860.9
$ws.Range("K116").Value = 6860.9
$ws.Range("M116").Value = -4566.9

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10594.739
$ws.Range("I3").Value = 6860.9
$ws.Range("K3").Value = 6860.9
$ws.Range("M3").Value = -6746.9
$ws.Range("H20").Value = 4001.3635
$ws.Range("I20").Value = 3890.8076
$ws.Range("J20").Value = 4412
$ws.Range("K20").Value = 3890.8076
$ws.Range("L20").Value = 4412
$ws.Range("M20").Value = -3643.8076
$ws.Range("N20").Value = -4906
$ws.Range("H138").Value = 53058.535
$ws.Range("J138").Value = 56496.92
$ws.Range("L138").Value = 56496.92
$ws.Range("N138").Value = -66776.92

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 359978.22
$ws.Range("I31").Value = 528346.1
$ws.Range("J31").Value = 4534.8887
$ws.Range("K31").Value = 528346.1
$ws.Range("L31").Value = 4534.8887
$ws.Range("M31").Value = -528051.1
$ws.Range("N31").Value = -5124.8887
$ws.Range("H34").Value = 359978.22
$ws.Range("I34").Value = 528346.1
$ws.Range("J34").Value = 4534.8887
$ws.Range("K34").Value = 528346.1
$ws.Range("L34").Value = 4534.8887
$ws.Range("M34").Value = -528144.1
$ws.Range("N34").Value = -4938.8887
$ws.Range("H58").Value = 3968.5
$ws.Range("J58").Value = 6998.6665
$ws.Range("L58").Value = 6998.6665
$ws.Range("N58").Value = -7404.6665
$ws.Range("H68").Value = 43384.31
$ws.Range("J68").Value = 43384.31
$ws.Range("L68").Value = 43384.31
$ws.Range("N68").Value = -44882.31
$ws.Range("H71").Value = 43384.31
$ws.Range("J71").Value = 43384.31
$ws.Range("L71").Value = 130152.93
$ws.Range("N71").Value = -137640.93
$ws.Range("H107").Value = 2766.0461
$ws.Range("I107").Value = 507.0811
$ws.Range("K107").Value = 507.0811
$ws.Range("M107").Value = 1412.9189
$ws.Range("H134").Value = 3498.0833
$ws.Range("I134").Value = 3498.0833
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10494.2499
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7959.249899999999
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 3968.5
$ws.Range("J136").Value = 6998.6665
$ws.Range("L136").Value = 20995.9995
$ws.Range("N136").Value = -26095.9995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 785.4
$ws.Range("J52").Value = 808.25
$ws.Range("L52").Value = 2424.75
$ws.Range("N52").Value = -2956.75
$ws.Range("H69").Value = 3700
$ws.Range("I69").Value = 1400
$ws.Range("J69").Value = 4083.3333
$ws.Range("K69").Value = 4200
$ws.Range("L69").Value = 12249.9999
$ws.Range("M69").Value = -3389
$ws.Range("N69").Value = -13871.9999
$ws.Range("H72").Value = 3700
$ws.Range("I72").Value = 1400
$ws.Range("J72").Value = 4083.3333
$ws.Range("K72").Value = 12600
$ws.Range("L72").Value = 36749.9997
$ws.Range("M72").Value = -8544
$ws.Range("N72").Value = -44861.9997
$ws.Range("H97").Value = 424.4
$ws.Range("J97").Value = 476.33334
$ws.Range("L97").Value = 1429.00002
$ws.Range("N97").Value = -2421.00002
$ws.Range("H131").Value = 16130734
$ws.Range("I131").Value = 83334536
$ws.Range("J131").Value = 1820.16
$ws.Range("K131").Value = 250003608
$ws.Range("L131").Value = 5460.48
$ws.Range("M131").Value = -249998568
$ws.Range("N131").Value = -15540.48

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 57000
$ws.Range("J86").Value = 57000
$ws.Range("L86").Value = 57000
$ws.Range("N86").Value = -59372
$ws.Range("H89").Value = 57000
$ws.Range("J89").Value = 57000
$ws.Range("L89").Value = 171000
$ws.Range("N89").Value = -182856
$ws.Range("H102").Value = 45823.918
$ws.Range("I102").Value = 3925.3684
$ws.Range("J102").Value = 205038.4
$ws.Range("K102").Value = 3925.3684
$ws.Range("L102").Value = 205038.4
$ws.Range("M102").Value = -2303.3684
$ws.Range("N102").Value = -208282.4
$ws.Range("H104").Value = 57160
$ws.Range("J104").Value = 57160
$ws.Range("L104").Value = 57160
$ws.Range("N104").Value = -64148
$ws.Range("H113").Value = 2696
$ws.Range("I113").Value = 2182
$ws.Range("K113").Value = 2182
$ws.Range("M113").Value = -12
$ws.Range("H126").Value = 9267.037
$ws.Range("I126").Value = 10144.826
$ws.Range("J126").Value = 4219.75
$ws.Range("K126").Value = 30434.478
$ws.Range("L126").Value = 12659.25
$ws.Range("M126").Value = -27964.478
$ws.Range("N126").Value = -17599.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 723.5
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 631.3333
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 631.3333
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1221.3333
$ws.Range("H27").Value = 723.5
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 631.3333
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 631.3333
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -845.3333
$ws.Range("H55").Value = 5196.273
$ws.Range("I55").Value = 661.75
$ws.Range("K55").Value = 661.75
$ws.Range("M55").Value = -488.75
$ws.Range("H82").Value = 1656.2413
$ws.Range("I82").Value = 1508.1305
$ws.Range("K82").Value = 1508.1305
$ws.Range("M82").Value = -1147.1305
$ws.Range("H85").Value = 1656.2413
$ws.Range("I85").Value = 1508.1305
$ws.Range("K85").Value = 1508.1305
$ws.Range("M85").Value = -260.1305
$ws.Range("H132").Value = 4500.407
$ws.Range("I132").Value = 4335.6313
$ws.Range("J132").Value = 4891.75
$ws.Range("K132").Value = 13006.8939
$ws.Range("L132").Value = 14675.25
$ws.Range("M132").Value = -10476.8939
$ws.Range("N132").Value = -19735.25
$ws.Range("H136").Value = 2476.8462
$ws.Range("I136").Value = 1790.8636
$ws.Range("K136").Value = 5372.5908
$ws.Range("M136").Value = -2822.5908
$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 20001
$ws.Range("I81").Value = 20001
$ws.Range("K81").Value = 40002
$ws.Range("M81").Value = -38941
$ws.Range("H84").Value = 20001
$ws.Range("I84").Value = 20001
$ws.Range("K84").Value = 200010
$ws.Range("M84").Value = -194706
$ws.Range("H96").Value = 4279.222
$ws.Range("I96").Value = 4002.5
$ws.Range("J96").Value = 4832.6665
$ws.Range("K96").Value = 4002.5
$ws.Range("L96").Value = 4832.6665
$ws.Range("M96").Value = -2629.5
$ws.Range("N96").Value = -7578.6665
$ws.Range("H126").Value = 2076.5715
$ws.Range("I126").Value = 1562.5385
$ws.Range("J126").Value = 2911.875
$ws.Range("K126").Value = 4687.6155
$ws.Range("L126").Value = 8735.625
$ws.Range("M126").Value = -2217.6155
$ws.Range("N126").Value = -13675.625
$ws.Range("H136").Value = 361332.1
$ws.Range("I136").Value = 458561.47
$ws.Range("K136").Value = 1375684.41
$ws.Range("M136").Value = -1373134.41
